$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 150.3
$ws.Cells.Item(12, 9).Value = 140.2
$ws.Cells.Item(12, 10).Value = 160.4
$ws.Cells.Item(12, 11).Value = 140.2
$ws.Cells.Item(12, 12).Value = 160.4
$ws.Cells.Item(12, 13).Value = 29.80000000000001
$ws.Cells.Item(12, 14).Value = -500.4

$ws.Cells.Item(15, 8).Value = 213466.2
$ws.Cells.Item(15, 9).Value = 213466.2
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 640398.6000000001
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -640229.6000000001

$ws.Cells.Item(19, 8).Value = 1223.1305
$ws.Cells.Item(19, 9).Value = 874
$ws.Cells.Item(19, 10).Value = 1677
$ws.Cells.Item(19, 11).Value = 874
$ws.Cells.Item(19, 12).Value = 1677
$ws.Cells.Item(19, 13).Value = -699
$ws.Cells.Item(19, 14).Value = -2027

$ws.Cells.Item(41, 8).Value = 482.0909
$ws.Cells.Item(41, 9).Value = 765
$ws.Cells.Item(41, 10).Value = 142.6
$ws.Cells.Item(41, 11).Value = 765
$ws.Cells.Item(41, 12).Value = 142.6
$ws.Cells.Item(41, 13).Value = -325
$ws.Cells.Item(41, 14).Value = -1022.6

$ws.Cells.Item(76, 8).Value = 3516.1892
$ws.Cells.Item(76, 9).Value = 3002.5862
$ws.Cells.Item(76, 10).Value = 5378
$ws.Cells.Item(76, 11).Value = 3002.5862
$ws.Cells.Item(76, 12).Value = 5378
$ws.Cells.Item(76, 13).Value = -2687.5862
$ws.Cells.Item(76, 14).Value = -6008

$ws.Cells.Item(79, 8).Value = 3516.1892
$ws.Cells.Item(79, 9).Value = 3002.5862
$ws.Cells.Item(79, 10).Value = 5378
$ws.Cells.Item(79, 11).Value = 3002.5862
$ws.Cells.Item(79, 12).Value = 5378
$ws.Cells.Item(79, 13).Value = -1910.5862
$ws.Cells.Item(79, 14).Value = -7562

$ws.Cells.Item(112, 8).Value = 1097
$ws.Cells.Item(112, 9).Value = 575
$ws.Cells.Item(112, 10).Value = 1174.3334
$ws.Cells.Item(112, 11).Value = 1725
$ws.Cells.Item(112, 12).Value = 3523.0002
$ws.Cells.Item(112, 13).Value = -617
$ws.Cells.Item(112, 14).Value = -5739.0002

$ws.Cells.Item(116, 8).Value = 2625
$ws.Cells.Item(116, 9).Value = 2833.3333
$ws.Cells.Item(116, 10).Value = 2000
$ws.Cells.Item(116, 11).Value = 2833.3333
$ws.Cells.Item(116, 12).Value = 2000
$ws.Cells.Item(116, 13).Value = 608.6667000000002
$ws.Cells.Item(116, 14).Value = -8884

$ws.Cells.Item(129, 8).Value = 994.1268
$ws.Cells.Item(129, 9).Value = 452.42856
$ws.Cells.Item(129, 10).Value = 1053.375
$ws.Cells.Item(129, 11).Value = 1357.28568
$ws.Cells.Item(129, 12).Value = 3160.125
$ws.Cells.Item(129, 13).Value = 3642.71432
$ws.Cells.Item(129, 14).Value = -13160.125

$ws.Cells.Item(135, 8).Value = 582.0909
$ws.Cells.Item(135, 9).Value = 566.9524
$ws.Cells.Item(135, 10).Value = 900
$ws.Cells.Item(135, 11).Value = 5102.5716
$ws.Cells.Item(135, 12).Value = 8100
$ws.Cells.Item(135, 13).Value = -2567.5716
$ws.Cells.Item(135, 14).Value = -13170

$ws.Cells.Item(137, 8).Value = 1392.2307
$ws.Cells.Item(137, 9).Value = 1091.3334
$ws.Cells.Item(137, 10).Value = 5003
$ws.Cells.Item(137, 11).Value = 3274.0002
$ws.Cells.Item(137, 12).Value = 15009
$ws.Cells.Item(137, 13).Value = -724.0001999999999
$ws.Cells.Item(137, 14).Value = -20109

$ws.Cells.Item(138, 8).Value = 1646.9661
$ws.Cells.Item(138, 9).Value = 1232.4584
$ws.Cells.Item(138, 10).Value = 3455.7273
$ws.Cells.Item(138, 11).Value = 3697.3752
$ws.Cells.Item(138, 12).Value = 10367.1819
$ws.Cells.Item(138, 13).Value = 1442.6248
$ws.Cells.Item(138, 14).Value = -20647.1819

$ws.Cells.Item(141, 8).Value = 8293.152
$ws.Cells.Item(141, 9).Value = 1285.5518
$ws.Cells.Item(141, 10).Value = 20247.295
$ws.Cells.Item(141, 11).Value = 3856.6554
$ws.Cells.Item(141, 12).Value = 60741.88499999999
$ws.Cells.Item(141, 13).Value = 1323.3446
$ws.Cells.Item(141, 14).Value = -71101.88499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 30301.2
$ws.Cells.Item(6, 9).Value = 37626.5
$ws.Cells.Item(6, 10).Value = 1000
$ws.Cells.Item(6, 11).Value = 37626.5
$ws.Cells.Item(6, 12).Value = 1000
$ws.Cells.Item(6, 13).Value = -37453.5
$ws.Cells.Item(6, 14).Value = -1346

$ws.Cells.Item(61, 8).Value = 1000
$ws.Cells.Item(61, 9).Value = 1000
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1000
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = -788

$ws.Cells.Item(136, 8).Value = 1000
$ws.Cells.Item(136, 9).Value = 1000
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 3000
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).Value = -450

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3371.25
$ws.Cells.Item(86, 9).Value = 3628.3333
$ws.Cells.Item(86, 10).Value = 2600
$ws.Cells.Item(86, 11).Value = 3628.3333
$ws.Cells.Item(86, 12).Value = 2600
$ws.Cells.Item(86, 13).Value = -2505.3333
$ws.Cells.Item(86, 14).Value = -4846

$ws.Cells.Item(89, 8).Value = 3371.25
$ws.Cells.Item(89, 9).Value = 3628.3333
$ws.Cells.Item(89, 10).Value = 2600
$ws.Cells.Item(89, 11).Value = 18141.6665
$ws.Cells.Item(89, 12).Value = 13000
$ws.Cells.Item(89, 13).Value = -12525.6665
$ws.Cells.Item(89, 14).Value = -24232

$ws.Cells.Item(128, 8).Value = 2000
$ws.Cells.Item(128, 9).Value = 2000
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 6000
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 13).Value = -3510

$ws.Cells.Item(134, 8).Value = 2355.5
$ws.Cells.Item(134, 9).Value = 2134
$ws.Cells.Item(134, 10).Value = 2665.6
$ws.Cells.Item(134, 11).Value = 6402
$ws.Cells.Item(134, 12).Value = 7996.799999999999
$ws.Cells.Item(134, 13).Value = -3867
$ws.Cells.Item(134, 14).Value = -13066.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2919
$ws.Cells.Item(86, 9).Value = 3083.75
$ws.Cells.Item(86, 10).Value = 2739.2727
$ws.Cells.Item(86, 11).Value = 3083.75
$ws.Cells.Item(86, 12).Value = 2739.2727
$ws.Cells.Item(86, 13).Value = -1960.75
$ws.Cells.Item(86, 14).Value = -4985.2727

$ws.Cells.Item(89, 8).Value = 2919
$ws.Cells.Item(89, 9).Value = 3083.75
$ws.Cells.Item(89, 10).Value = 2739.2727
$ws.Cells.Item(89, 11).Value = 15418.75
$ws.Cells.Item(89, 12).Value = 13696.3635
$ws.Cells.Item(89, 13).Value = -9802.75
$ws.Cells.Item(89, 14).Value = -24928.3635

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 1220
$ws.Cells.Item(22, 9).Value = 800
$ws.Cells.Item(22, 10).Value = 1966.6666
$ws.Cells.Item(22, 11).Value = 2400
$ws.Cells.Item(22, 12).Value = 5899.9998
$ws.Cells.Item(22, 13).Value = -2231
$ws.Cells.Item(22, 14).Value = -6237.9998

$ws.Cells.Item(27, 8).Value = 1220
$ws.Cells.Item(27, 9).Value = 800
$ws.Cells.Item(27, 10).Value = 1966.6666
$ws.Cells.Item(27, 11).Value = 2400
$ws.Cells.Item(27, 12).Value = 5899.9998
$ws.Cells.Item(27, 13).Value = -2298
$ws.Cells.Item(27, 14).Value = -6103.9998

$ws.Cells.Item(122, 8).Value = 668
$ws.Cells.Item(122, 9).Value = 434
$ws.Cells.Item(122, 10).Value = 1448
$ws.Cells.Item(122, 11).Value = 3906
$ws.Cells.Item(122, 12).Value = 13032
$ws.Cells.Item(122, 13).Value = -1456
$ws.Cells.Item(122, 14).Value = -17932

$ws.Cells.Item(131, 8).Value = 4302.2188
$ws.Cells.Item(131, 9).Value = 548.2857
$ws.Cells.Item(131, 10).Value = 5353.32
$ws.Cells.Item(131, 11).Value = 1644.8571
$ws.Cells.Item(131, 12).Value = 16059.96
$ws.Cells.Item(131, 13).Value = 3395.1429
$ws.Cells.Item(131, 14).Value = -26139.96

$ws.Cells.Item(132, 8).Value = 1764.8334
$ws.Cells.Item(132, 9).Value = 1648.25
$ws.Cells.Item(132, 10).Value = 1823.125
$ws.Cells.Item(132, 11).Value = 14834.25
$ws.Cells.Item(132, 12).Value = 16408.125
$ws.Cells.Item(132, 13).Value = -12304.25
$ws.Cells.Item(132, 14).Value = -21468.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4971.8823
$ws.Cells.Item(70, 9).Value = 4368.972
$ws.Cells.Item(70, 10).Value = 6418.8667
$ws.Cells.Item(70, 11).Value = 4368.972
$ws.Cells.Item(70, 12).Value = 6418.8667
$ws.Cells.Item(70, 13).Value = -4098.972
$ws.Cells.Item(70, 14).Value = -6958.8667

$ws.Cells.Item(73, 8).Value = 4971.8823
$ws.Cells.Item(73, 9).Value = 4368.972
$ws.Cells.Item(73, 10).Value = 6418.8667
$ws.Cells.Item(73, 11).Value = 4368.972
$ws.Cells.Item(73, 12).Value = 6418.8667
$ws.Cells.Item(73, 13).Value = -3432.972
$ws.Cells.Item(73, 14).Value = -8290.866699999999

$ws.Cells.Item(80, 8).Value = 2785.7144
$ws.Cells.Item(80, 9).Value = 2111.111
$ws.Cells.Item(80, 10).Value = 4000
$ws.Cells.Item(80, 11).Value = 2111.111
$ws.Cells.Item(80, 12).Value = 4000
$ws.Cells.Item(80, 13).Value = -1113.111
$ws.Cells.Item(80, 14).Value = -5996

$ws.Cells.Item(83, 8).Value = 2785.7144
$ws.Cells.Item(83, 9).Value = 2111.111
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 11).Value = 10555.555
$ws.Cells.Item(83, 12).Value = 20000
$ws.Cells.Item(83, 13).Value = -5563.555
$ws.Cells.Item(83, 14).Value = -29984

$ws.Cells.Item(113, 8).Value = 1891.2
$ws.Cells.Item(113, 9).Value = 1757.9286
$ws.Cells.Item(113, 10).Value = 2202.1667
$ws.Cells.Item(113, 11).Value = 1757.9286
$ws.Cells.Item(113, 12).Value = 2202.1667
$ws.Cells.Item(113, 13).Value = 412.0714
$ws.Cells.Item(113, 14).Value = -6542.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 715.9583
$ws.Cells.Item(16, 9).Value = 399
$ws.Cells.Item(16, 10).Value = 1920.4
$ws.Cells.Item(16, 11).Value = 399
$ws.Cells.Item(16, 12).Value = 1920.4
$ws.Cells.Item(16, 13).Value = -229
$ws.Cells.Item(16, 14).Value = -2260.4

$ws.Cells.Item(22, 8).Value = 774.9091
$ws.Cells.Item(22, 9).Value = 404
$ws.Cells.Item(22, 10).Value = 1084
$ws.Cells.Item(22, 11).Value = 404
$ws.Cells.Item(22, 12).Value = 1084
$ws.Cells.Item(22, 13).Value = -109
$ws.Cells.Item(22, 14).Value = -1674

$ws.Cells.Item(27, 8).Value = 774.9091
$ws.Cells.Item(27, 9).Value = 404
$ws.Cells.Item(27, 10).Value = 1084
$ws.Cells.Item(27, 11).Value = 404
$ws.Cells.Item(27, 12).Value = 1084
$ws.Cells.Item(27, 13).Value = -297
$ws.Cells.Item(27, 14).Value = -1298

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1812.2106
$ws.Cells.Item(132, 9).Value = 1561.7037
$ws.Cells.Item(132, 10).Value = 2427.0908
$ws.Cells.Item(132, 11).Value = 4685.1111
$ws.Cells.Item(132, 12).Value = 7281.2724
$ws.Cells.Item(132, 13).Value = -2155.1111
$ws.Cells.Item(132, 14).Value = -12341.2724

$ws.Cells.Item(136, 8).Value = 1977.7368
$ws.Cells.Item(136, 9).Value = 931.5
$ws.Cells.Item(136, 10).Value = 3771.2856
$ws.Cells.Item(136, 11).Value = 2794.5
$ws.Cells.Item(136, 12).Value = 11313.8568
$ws.Cells.Item(136, 13).Value = -244.5
$ws.Cells.Item(136, 14).Value = -16413.8568
